# Auto-upload VRF Excel file
# Adds a new worksheet named "asd" at the end of the workbook with the
# standard VRF outdoor/indoor model header row (matching the format used
# by every other sheet in this workbook).

$wb = $excel.ActiveWorkbook

# Insert the new sheet after the last existing sheet so it lands at the end.
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "asd"

# Header row values.
$ws.Cells.Item(1, 1).Value = "Outdoor Model"
$ws.Cells.Item(1, 2).Value = "Outdoor Quantity"
$ws.Cells.Item(1, 3).Value = "Outdoor Serial(s)"
$ws.Cells.Item(1, 4).Value = "Indoor Model"
$ws.Cells.Item(1, 5).Value = "Indoor Quantity"
$ws.Cells.Item(1, 6).Value = "Indoor Serial(s)"

# Match the bold / centered / thin-bordered header formatting used on
# every other sheet in the workbook.
$headerRange = $ws.Range("A1:F1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1
$headerRange.Borders.Weight = 2
